$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")
$ws.Activate()

# --- Row 29 : new daily-tracker entry (task #20) ---------------------------
$ws.Range("A29").Value = 20
$ws.Range("B29").Value = 44616              # 24-Feb-2022 (date serial, cell already formatted)
$ws.Range("C29").Value = "RPA GSS"
$ws.Range("D29").Value = "1. All 11 tasks of monthly has been implemented with MFA with OTP reading, and upload status, tested and it is running smoothly"
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Completed"

# --- Row 30 : second comment line that belongs to task #20 -----------------
$ws.Range("D30").Value = "2. Whereas , the task of invoice generation is work in progress"

$e30 = $ws.Range("E30")
$e30.NumberFormat = "0%"
$e30.Value = 0.6

$f30 = $ws.Range("F30")
$f30.Value = "WIP"
$f30.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$f30.Borders.Item(7).Weight = 2      # xlThin
$f30.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$f30.Borders.Item(10).Weight = 2     # xlThin

# --- keep the selection where the user left it off --------------------------
$ws.Range("F30").Select()
